# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker data table (rows 16-29) is reorganized: instead of being grouped
# by worker (all periods for ILMER, then all periods for IVANNA), it is now
# grouped by period (both workers for 2103, then both workers for 2104, etc.),
# and a new worker (IVANNA PAOLA PASSO CORREA / 1007254953) is added for
# period 2103. The "Valor Mora" (column G) is updated to 877803 for every
# row, and "Salario Basico" (column F) becomes 35112 for every row except
# the most recent period (2109), which keeps 23408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docType = "CC"
$idIlmer = "12598379"
$nameIlmer = "ILMER IVAN PASSO PUELLO"
$idIvanna = "1007254953"
$nameIvanna = "IVANNA PAOLA PASSO CORREA"

$periods = @("2103", "2104", "2105", "2106", "2107", "2108", "2109")

$row = 16
foreach ($periodo in $periods) {
    if ($periodo -eq "2109") {
        $salario = 23408
    } else {
        $salario = 35112
    }
    $mora = 877803

    # ILMER IVAN PASSO PUELLO row
    $ws.Cells.Item($row, 2).Value = $docType
    $ws.Cells.Item($row, 3).Value = $idIlmer
    $ws.Cells.Item($row, 4).Value = $nameIlmer
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $salario
    $ws.Cells.Item($row, 7).Value = $mora
    $row = $row + 1

    # IVANNA PAOLA PASSO CORREA row
    $ws.Cells.Item($row, 2).Value = $docType
    $ws.Cells.Item($row, 3).Value = $idIvanna
    $ws.Cells.Item($row, 4).Value = $nameIvanna
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $salario
    $ws.Cells.Item($row, 7).Value = $mora
    $row = $row + 1
}
